$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A to hold the "Avg_Quantity" values
$ws.Columns("A:A").Insert()

# Header row
$ws.Range("A1").Value = "Avg_Quantity"
$ws.Range("B1").Value = "City"
$ws.Range("C1").Value = "Quantity"

# Remove the old "Timbre" column (now column D after insertion)
$ws.Columns("D:D").Delete()

# Match A1's formatting to the other header cells (bold, bordered)
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# Data rows
$ws.Range("A2").Value = 77.75
$ws.Range("B2").Value = "ECOUFLANT"
$ws.Range("C2").Value = 933

$ws.Range("A3").Value = 66.78571428571429
$ws.Range("B3").Value = "MONTIGNE LE BRILLANT"
$ws.Range("C3").Value = 935

$ws.Range("A4").Value = 92
$ws.Range("B4").Value = "LES GARENNES SUR LOIRE"
$ws.Range("C4").Value = 1196

$ws.Range("A5").Value = 64.1025641025641
$ws.Range("B5").Value = "RENAZE"
$ws.Range("C5").Value = 2500

$ws.Range("A6").Value = 90.29166666666667
$ws.Range("B6").Value = "TORCE VIVIERS EN CHARNIE"
$ws.Range("C6").Value = 2167
